$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '43.739.14'
$ws.Range("E2").Value = '  +0.86%  '
Set-TextValue $ws.Range("D3") '2.353.73'
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("E4").Value = '  -0.11%  '
Set-TextValue $ws.Range("D5") '0.672'
$ws.Range("E5").Value = '  +3.76%  '
Set-TextValue $ws.Range("D6") '235.55'
$ws.Range("E6").Value = '  +1.33%  '
Set-TextValue $ws.Range("D7") '73.84'
$ws.Range("E7").Value = '  +10.96%  '
$ws.Range("E8").Value = '  -0.01%  '
Set-TextValue $ws.Range("D9") '0.543'
$ws.Range("E9").Value = '  +19.47%  '
Set-TextValue $ws.Range("D10") '0.0988'
$ws.Range("E10").Value = '  +2.06%  '
Set-TextValue $ws.Range("D11") '28.31'
$ws.Range("E11").Value = '  +5.61%  '
$ws.Range("E12").Value = '  +1.95%  '
Set-TextValue $ws.Range("D13") '2.702.41'
$ws.Range("E13").Value = '  +0.80%  '
Set-TextValue $ws.Range("D14") '16.75'
$ws.Range("E14").Value = '  +7.95%  '
Set-TextValue $ws.Range("D15") '6.70'
$ws.Range("E15").Value = '  +6.99%  '
Set-TextValue $ws.Range("D16") '0.891'
$ws.Range("E16").Value = '  +4.31%  '
Set-TextValue $ws.Range("D17") '2.414.26'
$ws.Range("E17").Value = '  +4.13%  '
Set-TextValue $ws.Range("D18") '43.770.67'
$ws.Range("E18").Value = '  +1.15%  '
Set-TextValue $ws.Range("D19") '0.0000102'
$ws.Range("E19").Value = '  +3.56%  '
$ws.Range("E20").Value = '  +3.97%  '
Set-TextValue $ws.Range("D21") '6.44'
$ws.Range("E21").Value = '  +3.05%  '
Set-TextValue $ws.Range("D22") '253.66'
$ws.Range("E22").Value = '  +2.10%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  -1.82%  '
$ws.Range("E25").Value = '  +2.65%  '
Set-TextValue $ws.Range("D26") '10.58'
$ws.Range("E26").Value = '  +6.52%  '
Set-TextValue $ws.Range("D27") '2.29'
$ws.Range("E27").Value = '  +0.47%  '
Set-TextValue $ws.Range("D28") '22.38'
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D29") '172.75'
$ws.Range("E29").Value = '  -0.81%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D30") '1.59'
$ws.Range("E30").Value = '  +7.54%  '
$ws.Range("E31").Value = '  +1.86%  '
$ws.Range("E32").Value = '  +5.21%  '
$ws.Range("E33").Value = '  +2.92%  '
$ws.Range("E34").Value = '  +3.47%  '
Set-TextValue $ws.Range("D35") '5.15'
$ws.Range("E35").Value = '  +3.70%  '
$ws.Range("E36").Value = '  +7.24%  '
$ws.Range("E37").Value = '  -4.05%  '
Set-TextValue $ws.Range("D38") '6.40'
$ws.Range("E38").Value = '  -1.59%  '
$ws.Range("E39").Value = '  +5.91%  '
Set-TextValue $ws.Range("D40") '19.51'
$ws.Range("E40").Value = '  +6.21%  '
$ws.Range("E41").Value = '  -0.14%  '
Set-TextValue $ws.Range("D42") '8.86'
$ws.Range("E42").Value = '  -2.63%  '
$ws.Range("E43").Value = '  +1.34%  '
$ws.Range("E44").Value = '  +3.28%  '
$ws.Range("E45").Value = '  -0.76%  '
Set-TextValue $ws.Range("D46") '4.44'
$ws.Range("E46").Value = '  +2.16%  '
Set-TextValue $ws.Range("D47") '0.181'
$ws.Range("E47").Value = '  +11.28%  '
Set-TextValue $ws.Range("D48") '97.13'
$ws.Range("E48").Value = '  -2.30%  '
Set-TextValue $ws.Range("D49") '1.437.67'
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("E50").Value = '  +1.43%  '
Set-TextValue $ws.Range("D51") '2.577.04'
$ws.Range("E51").Value = '  +0.90%  '
